# Team builder results - reshuffle team rosters and GPA values
# (matches the diff: shared-string order churn from inserted/reordered team
#  rosters, plus new GPA numbers in columns H:K/L for every data row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
  2  = @{ "B"="EE08";  "C"="ME03";  "D"="ME05";  "E"="ME25";                  "H"=3.2; "I"=3.3; "J"=2.5; "K"=3.5 }
  3  = @{ "B"="EE09";  "C"="ME24";  "D"="ME28";  "E"="ME45";                  "H"=2.3; "I"=3.4; "J"=3.8; "K"=3.5 }
  4  = @{ "B"="EE02";  "C"="EE03";  "D"="ME42";  "E"="ME48";                  "H"=3.6; "I"=3.7; "J"=2.2; "K"=2.8 }
  5  = @{ "B"="EE07";  "C"="ME35";  "D"="ME53";  "E"="ME57";                  "H"=3.1; "I"=2.5; "J"=3.3; "K"=3.7 }
  6  = @{ "B"="EE17";  "C"="ME19";  "D"="ME32";                               "H"=3.1; "I"=3.9; "J"=2.2           }
  7  = @{ "B"="EE00";  "C"="ME39";  "D"="ME50";  "E"="ME60";                  "H"=3.4; "I"=3.9; "J"=3.0; "K"=3.0 }
  8  = @{ "B"="EE23";  "C"="ME09";  "D"="ME23";  "E"="ME40";                  "H"=3.7; "I"=2.9; "J"=3.3; "K"=3.0 }
  9  = @{ "B"="EE04";  "C"="ME06";  "D"="ME30";  "E"="ME61";                  "H"=3.8; "I"=3.6; "J"=3.0; "K"=2.1 }
  10 = @{ "B"="EE06";  "C"="ME37";  "D"="ME52";  "E"="ME63";                  "H"=3.0; "I"=3.7; "J"=3.2; "K"=3.3 }
  11 = @{ "B"="EE24";  "C"="ME07";  "D"="ME54";  "E"="ME58";                  "H"=3.8; "I"=3.7; "J"=2.4; "K"=3.8 }
  12 = @{ "B"="EE10";  "C"="ME04";  "D"="ME10";  "E"="ME62";                  "H"=3.4; "I"=3.4; "J"=3.0; "K"=3.2 }
  13 = @{ "B"="EE11";  "C"="ME11";  "D"="ME27";  "E"="ME34";                  "H"=3.5; "I"=3.1; "J"=2.7; "K"=3.4 }
  14 = @{ "B"="EE12";  "C"="ME12";  "D"="ME33";  "E"="ME44";                  "H"=3.6; "I"=3.2; "J"=3.3; "K"=2.4 }
  15 = @{ "B"="CpE01"; "C"="CpE03"; "D"="EE13";  "E"="ME13";                  "H"=3.7; "I"=3.9; "J"=2.7; "K"=3.3 }
  16 = @{ "B"="EE14";  "C"="ME14";  "D"="ME46";  "E"="ME47";                  "H"=3.8; "I"=2.4; "J"=3.6; "K"=3.7 }
  17 = @{ "B"="CpE00"; "C"="EE15";  "D"="ME15";  "E"="ME36";                  "H"=2.6; "I"=3.9; "J"=3.5; "K"=3.6 }
  18 = @{ "B"="EE16";  "C"="ME16";  "D"="ME26";  "E"="ME41";                  "H"=3.0; "I"=3.6; "J"=3.6; "K"=3.1 }
  19 = @{ "B"="EE01";  "C"="ME00";  "D"="ME17";  "E"="ME49"; "F"="ME55";      "H"=2.5; "I"=2.0; "J"=3.7; "K"=3.9; "L"=3.5 }
  20 = @{ "B"="EE18";  "C"="ME02";  "D"="ME18";  "E"="ME59";                  "H"=3.2; "I"=2.2; "J"=3.8; "K"=3.9 }
  21 = @{ "B"="CpE02"; "C"="EE19";  "D"="ME08";  "E"="ME29";                  "H"=3.8; "I"=2.3; "J"=3.8; "K"=3.9 }
  22 = @{ "B"="EE20";  "C"="ME21";  "D"="ME31";  "E"="ME56";                  "H"=3.4; "I"=3.1; "J"=3.1; "K"=3.6 }
  23 = @{ "B"="EE21";  "C"="ME20";  "D"="ME38";  "E"="ME43";                  "H"=3.5; "I"=2.0; "J"=3.8; "K"=3.3 }
  24 = @{ "B"="EE05";  "C"="EE22";  "D"="ME01";  "E"="ME22"; "F"="ME51";      "H"=3.9; "I"=3.6; "J"=3.1; "K"=3.2; "L"=3.1 }
  25 = @{ "B"="CE00";  "C"="CE02";  "D"="CE03";  "E"="CE05";                  "H"=3.9; "I"=3.1; "J"=3.2; "K"=3.4 }
  26 = @{ "B"="CE01";  "C"="CE04";  "D"="CE06";                               "H"=2.7; "I"=3.3; "J"=3.5           }
}

foreach ($r in $rows.Keys) {
  $rowData = $rows[$r]
  foreach ($col in $rowData.Keys) {
    $ws.Cells.Item([int]$r, [int]([int][char]$col - [int][char]'A' + 1)).Value = $rowData[$col]
  }
}
